$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

# --- Column widths (narrower layout) ---
$ws.Columns.Item(1).ColumnWidth = 13.5
$ws.Columns.Item(2).ColumnWidth = 5.166666666666667
$ws.Columns.Item(3).ColumnWidth = 13.5
$ws.Columns.Item(4).ColumnWidth = 10.0
$ws.Columns.Item(5).ColumnWidth = 10.0
$ws.Columns.Item(6).ColumnWidth = 8.833333333333334
$ws.Columns.Item(7).ColumnWidth = 10.0
$ws.Columns.Item(8).ColumnWidth = 14.833333333333334

# --- Header row: insert line breaks between words ---
$ws.Range("B1").Value = "Kia" + $nl + "Rio"
$ws.Range("C1").Value = "Volkswagen" + $nl + "Golf"
$ws.Range("D1").Value = "Toyota" + $nl + "Corolla"
$ws.Range("E1").Value = "Skoda" + $nl + "Octavia"
$ws.Range("F1").Value = "BMW" + $nl + "3" + $nl + "Series"
$ws.Range("G1").Value = "Hyundai" + $nl + "Solaris"
$ws.Range("H1").Value = "Вектор" + $nl + "приоритетов"

# --- Row labels (column A): same car names with line breaks ---
$ws.Range("A2").Value = "Kia" + $nl + "Rio"
$ws.Range("A3").Value = "Volkswagen" + $nl + "Golf"
$ws.Range("A4").Value = "Toyota" + $nl + "Corolla"
$ws.Range("A5").Value = "Skoda" + $nl + "Octavia"
$ws.Range("A6").Value = "BMW" + $nl + "3" + $nl + "Series"
$ws.Range("A7").Value = "Hyundai" + $nl + "Solaris"

# --- Priority vector value corrections ---
$ws.Range("H3").Value = 0.228
$ws.Range("H6").Value = 0.362
